$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# 1. Insert a new row at sheet row 22 (pushes existing data down by one row)
$ws.Rows.Item(22).Insert()

# 2. Grow the table (ListObject) so the new row is included in Table3
$lo.Resize($ws.Range("A1:K23"))

# 3. Populate the new row (row 22) with the "Protoboard" part purchased from Frys
$ws.Range("A22").Value = "Protoboard"
$ws.Range("B22").Value = "4x5"
$ws.Range("D22").Value = "Twin"
$ws.Range("E22").Value = "8000-45-LF"
$ws.Range("F22").Value = "Frys"
$ws.Range("G22").Value = 4986181
$ws.Range("H22").Value = "http://www.frys.com/product/4986181?source=googleps&gclid=CNPX5OnlxbQCFcxAMgodfm0AYw"
$ws.Range("I22").Value = 12.69
$ws.Range("J22").Value = 2
$ws.Range("K22").Formula = "=Table3[[#This Row],[Price]]*Table3[[#This Row],[Quantity]]"

# 4. Apply the "Bad" (red) cell style to the whole new row, matching the other
#    "needs attention" rows added to this BOM.
$ws.Range("A22:J22").Style = "Bad"
$ws.Range("C22").Style = "Bad"
$ws.Range("G22").Style = "Bad"
$ws.Range("G22").HorizontalAlignment = -4131
$ws.Range("I22").Style = "Bad"
$ws.Range("I22").NumberFormat = """$""#,##0.00_);[Red]\(""$""#,##0.00\)"

# 5. Hyperlink the distributor part number to the Frys product page
$ws.Hyperlinks.Add($ws.Range("G22"), "http://www.frys.com/product/4986181?source=googleps&gclid=CNPX5OnlxbQCFcxAMgodfm0AYw", "", "", "http://www.frys.com/product/4986181?source=googleps&gclid=CNPX5OnlxbQCFcxAMgodfm0AYw") | Out-Null

# 6. Keep the sheet's used range consistent with the extra row at the bottom
$ws.Rows.Item(186).EntireRow.RowHeight = $ws.Rows.Item(185).RowHeight

# 7. Restore the active selection like the author left it
$ws.Range("J22").Select()
